# Update the Jogos da Semana FlashScore sheet:
#  - tweak a handful of existing odds on row 3 and one value on row 5
#  - insert a brand-new match as row 8 (Pereira vs Atl. Nacional),
#    pushing the former row 8 (Tepatitlan de Morelos vs Tapatio) down to row 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a handful of existing odds on row 3 ---
$ws.Range("G3").Value  = 2.7
$ws.Range("I3").Value  = 2.9
$ws.Range("W3").Value  = 7.5
$ws.Range("X3").Value  = 12
$ws.Range("Z3").Value  = 26
$ws.Range("AB3").Value = 34
$ws.Range("AL3").Value = 26
$ws.Range("AW3").Value = 4.75

# --- Update a single value on row 5 ---
$ws.Range("BC5").Value = 126

# --- Insert a new row above row 8 (old row 8 shifts down to row 9) ---
$ws.Rows.Item(8).Insert()

# --- Populate the newly inserted row 8 with the new match data ---
$newRow = New-Object 'object[,]' 1,56
$newRow[0,0]  = "djGN4nIr"
$newRow[0,1]  = "PLACEHOLDER_DATE"
$newRow[0,2]  = "22:45"
$newRow[0,3]  = "COLOMBIA - PRIMERA A"
$newRow[0,4]  = "Pereira"
$newRow[0,5]  = "Atl. Nacional"
$newRow[0,6]  = 2.55
$newRow[0,7]  = 3
$newRow[0,8]  = 3
$newRow[0,9]  = 3.4
$newRow[0,10] = 1.91
$newRow[0,11] = 3.75
$newRow[0,12] = 1.1
$newRow[0,13] = 7
$newRow[0,14] = 1.5
$newRow[0,15] = 2.5
$newRow[0,16] = 2.5
$newRow[0,17] = 1.5
$newRow[0,18] = 1.57
$newRow[0,19] = 2.25
$newRow[0,20] = 2.1
$newRow[0,21] = 1.67
$newRow[0,22] = 6.5
$newRow[0,23] = 11
$newRow[0,24] = 11
$newRow[0,25] = 26
$newRow[0,26] = 26
$newRow[0,27] = 41
$newRow[0,28] = 6.5
$newRow[0,29] = 6
$newRow[0,30] = 19
$newRow[0,31] = 67
$newRow[0,32] = 900
$newRow[0,33] = 7
$newRow[0,34] = 13
$newRow[0,35] = 12
$newRow[0,36] = 34
$newRow[0,37] = 29
$newRow[0,38] = 41
$newRow[0,39] = 4.33
$newRow[0,40] = 15
$newRow[0,41] = 29
$newRow[0,42] = 51
$newRow[0,43] = 81
$newRow[0,44] = 301
$newRow[0,45] = 2.25
$newRow[0,46] = 9
$newRow[0,47] = 67
$newRow[0,48] = 4.75
$newRow[0,49] = 19
$newRow[0,50] = 34
$newRow[0,51] = 67
$newRow[0,52] = 101
$newRow[0,53] = 301
$newRow[0,54] = 126
$newRow[0,55] = 126

$ws.Range("A8:BD8").Value = $newRow

# The Date column (B) holds a dd/mm/yyyy-looking string that Excel would
# otherwise auto-convert into a date serial number. Force it to be written
# as plain text (matching the rest of the Date column), then drop the
# number-format override so the cell keeps the workbook's default style.
$dateCell = $ws.Range("B8")
$dateCell.NumberFormat = "@"
$dateCell.Value = "07/11/2024"
$dateCell.ClearFormats()
